$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are formatted as Text so Excel does not auto-convert
# numeric-looking / percent-looking strings into numbers when we set .Value
$updates = @{
    'D2' = '301.48'
    'E2' = '1.07%'
    'D3' = '31.82'
    'E3' = '1.82%'
    'D4' = '5.126'
    'E4' = '0.54%'
    'D5' = '0.07842'
    'D6' = '2.258'
    'E6' = '-10.22%'
    'D7' = '7.809'
    'D8' = '3.809'
    'E8' = '-0.04%'
    'D9' = '0.9282'
    'E9' = '1.08%'
    'D10' = '0.1774'
    'E10' = '2.75%'
    'D11' = '0.07691'
    'E11' = '4.83%'
    'D12' = '0.08911'
    'E12' = '3.92%'
    'D13' = '0.03101'
    'E13' = '2.48%'
    'E14' = '0.51%'
    'D15' = '0.001508'
    'E15' = '-0.02%'
    'D16' = '0.005804'
    'E16' = '-2.88%'
    'D17' = '3.480'
    'E17' = '-0.86%'
    'D18' = '2.250'
    'E18' = '0.17%'
    'E20' = '-0.77%'
    'D21' = '4.322'
    'E21' = '-6.10%'
    'D23' = '0.04599'
    'E23' = '-0.45%'
    'E24' = '0.41%'
    'D25' = '0.004485'
    'E25' = '1.06%'
    'E26' = '4.11%'
    'E27' = '-1.31%'
    'D39' = '0.01783'
    'D40' = '0.04793'
    'E40' = '6.14%'
    'D41' = '0.007354'
    'E41' = '3.97%'
    'D42' = '0.1365'
    'E42' = '1.73%'
    'D43' = '0.002189'
    'E43' = '-2.28%'
    'D44' = '0.009834'
    'E44' = '0.02%'
    'D45' = '0.00006264'
    'E45' = '-5.10%'
    'D46' = '0.00000000749'
    'E46' = '-0.05%'
    'D48' = '1.074'
    'E48' = '30.85%'
    'E49' = '-0.05%'
    'E50' = '-0.05%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

